$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in / update column F ("SolderPointsDIP") values for the refreshed
# KiCad BOM export.
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 17
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 164
$ws.Range("F12").Value = 211
$ws.Range("F13").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 47
$ws.Range("F22").Value = 15

# Update the selected cell to match the author's final cursor position.
$ws.Range("F23").Select()
